$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report's "totals" block (rows 8-12) and "main table" block (rows 2-6)
# swap places. Row 7 stays put as the blank separator between them.
#
# Move totals block (rows 8-12) out of the way to scratch rows first,
# then move main block (rows 2-6) down into rows 8-12, then move the
# totals block from scratch into rows 2-6.
$ws.Rows("8:12").Cut($ws.Rows("50:54"))
$ws.Rows("2:6").Cut($ws.Rows("8:12"))
$ws.Rows("50:54").Cut($ws.Rows("2:6"))

# Fix up the two header/title cells whose text changed.
$ws.Range("B8").Value = "Развернутый анализ продаж и прибыли"
$ws.Range("E2").Value = "Итоги анализа продаж и прибыли"

# Rebuild merged cells for the new layout.
$ws.Range("E3:F3").UnMerge()
$ws.Range("E5:F5").UnMerge()
$ws.Range("E4:F4").UnMerge()
$ws.Range("E2:J2").UnMerge()
$ws.Range("B8:J8").UnMerge()
$ws.Range("B9:E9").UnMerge()
$ws.Range("H9:I9").UnMerge()
$ws.Range("F9:G9").UnMerge()
$ws.Range("J9:J10").UnMerge()

$ws.Range("E3:F3").Merge()
$ws.Range("E5:F5").Merge()
$ws.Range("E4:F4").Merge()
$ws.Range("E2:J2").Merge()
$ws.Range("B8:J8").Merge()
$ws.Range("B9:E9").Merge()
$ws.Range("H9:I9").Merge()
$ws.Range("F9:G9").Merge()
$ws.Range("J9:J10").Merge()

# Update the saved selection.
$ws.Range("H16").Select()
